$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New accelerometer/gyroscope readings for the existing "falling" samples (rows 2-31,
# columns C:H) -- rows 2-21 are overwritten in place, rows 22-31 are brand new.
$chCsv = @"
-3.729709470272064,9.457800364494323,0.187229474633932,0.0209221355617046,-0.0198531206697225,-0.0239764600992202
-3.684299826622009,9.67076587677002,0.04418391920626155,0.06856962293386459,-0.06704246252775189,0.1534798890352249
-3.839793586730957,9.701957559585571,0.08396954610943785,0.0235183127224445,-0.0152716310694813,0.030695978552103
-3.847879505157471,9.791793012619019,0.2833944983780385,0.078801617026329,0.009010262787342,-0.0384845100343227
-3.597369003295899,9.792463493347167,0.5209171175956726,0.0911716371774673,0.0236710291355848,-0.0128281703218817
-3.581359481811524,9.727484250068665,0.4699014097452164,0.0241291765123605,-0.0087048299610614,-0.0103847095742821
-3.805133938789368,9.641714334487917,0.3316636346280575,0.0088575463742017,0.0404698215425014,0.0618501044809818
-3.644986343383789,9.586582708358765,0.1923969518393278,-0.0404698215425014,-0.0221438650041818,0.0143553335219621
-3.727581739425659,9.62705430984497,0.3099611997604369,0.0478002056479454,-0.022754730656743,0.08124507963657369
-3.772140645980835,9.524589729309081,0.2081888042390346,0.0325285755097866,0.0065668015740811,0.0574213340878486
-4.169818592071533,9.43976936340332,0.4335370913147927,-0.0219911485910415,0.0610865242779254,0.0355829000473022
-3.96463143825531,9.59680449962616,0.3559434115886688,-0.0038179077673703,0.0641408488154411,-0.0708603709936142
-3.995069694519043,9.550227546691893,0.2403170883655548,0.0064140851609408,-0.0274889357388019,-0.0403171069920063
-3.729709470272064,9.457800364494323,0.187229474633932,0.0209221355617046,-0.0198531206697225,-0.0239764600992202
-3.674985194206238,9.404587173461914,0.3419130772352217,-0.0395535230636596,0.0210748501121997,-0.0496328026056289
-3.945914745330811,9.500829362869263,0.6805746570229532,0.0119118718430399,0.0239764600992202,-0.0226020142436027
-3.8448965549469,9.112732887268068,1.50640323758125,0.0048869219608604,-0.0102319931611418,0.0267253536731004
-3.954125237464904,8.554272603988649,2.520113927125929,-0.1745547503232956,0.1945605874061584,0.0557414554059505
-4.40533800125122,8.189027309417723,3.551703810691835,-0.1134682223200798,0.6869179606437683,-0.2023491114377975
-4.528182744979861,7.154003858566278,4.620719850063325,-0.1960877478122711,0.6270532011985779,-0.06856962293386459
-5.236727142333985,5.354243278503416,5.061923789978028,0.6217080950737,0.2063197344541549,-1.041983366012573
-5.258063554763795,5.088717103004436,5.500368356704722,0.1171334087848663,1.447597861289978,-3.771023988723755
-4.872933006286623,10.60090007781982,2.667490434646607,-0.1485929638147354,-2.555707454681396,-3.479946613311768
-0.9181296348571681,8.491318988800044,0.2766812086105288,1.612989664077759,3.036305665969849,-3.449708700180054
7.330888652801525,7.476689434051512,3.909529781341558,-0.2490803003311157,-4.406629085540772,2.534785270690918
17.60857832431794,8.065273070335389,5.904095327854156,0.09239336848258969,-2.676048040390014,-0.3824016451835632
8.862465858459505,4.974191069602975,2.124806880950939,0.0197004042565822,-0.9019425511360168,0.6754642724990845
5.792156529426575,5.206499457359316,1.24618867635727,-0.173791155219078,-0.06688974797725671,-0.0691804885864257
6.776297664642335,7.347594261169435,2.1927396774292,-0.3798054754734039,0.1803579628467559,0.7594582438468933
7.413076007366181,7.669743299484253,3.059547257423401,0.0276416521519422,-0.9390525817871094,0.1821905523538589
"@
$chRows = $chCsv.Trim() -split "`n"
$chData = New-Object "object[,]" $chRows.Count,6
for ($i = 0; $i -lt $chRows.Count; $i++) {
    $parts = $chRows[$i].Split(",")
    for ($j = 0; $j -lt 6; $j++) {
        $chData[$i,$j] = [double]$parts[$j]
    }
}
$ws.Range("C2:H31").Value = $chData

# Timestamp + label columns (A:B) for the 10 newly appended rows (22-31), continuing
# the existing 100-unit timestamp sequence.
$abCsv = @"
2000,falling
2100,falling
2200,falling
2300,falling
2400,falling
2500,falling
2600,falling
2700,falling
2800,falling
2900,falling
"@
$abRows = $abCsv.Trim() -split "`n"
$abData = New-Object "object[,]" $abRows.Count,2
for ($i = 0; $i -lt $abRows.Count; $i++) {
    $parts = $abRows[$i].Split(",")
    $abData[$i,0] = [double]$parts[0]
    $abData[$i,1] = [string]$parts[1]
}
$ws.Range("A22:B31").Value = $abData
